# semana 49 de 2025
# Adds week 49 (column AZ) to the weekly IRA-Ext revision sheet and
# corrects several previously-reported weekly counts (rows 37, 43, 57)
# with updated/consolidated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header for week 49 (column AZ), stored as text like the other
#     week-number headers in row 1 (D1="1" ... AY1="48").
$ws.Cells.Item(1, 52).Value = "'49"

# --- New week-49 (column AZ) counts per UPGD row, plus the one newly
#     populated AY36 cell.
$ws.Cells.Item(2, 52).Value = 27
$ws.Cells.Item(5, 52).Value = 2
$ws.Cells.Item(6, 52).Value = 59
$ws.Cells.Item(7, 52).Value = 23
$ws.Cells.Item(8, 52).Value = 13
$ws.Cells.Item(9, 52).Value = 5
$ws.Cells.Item(10, 52).Value = 3
$ws.Cells.Item(11, 52).Value = 1
$ws.Cells.Item(12, 52).Value = 3
$ws.Cells.Item(14, 52).Value = 2
$ws.Cells.Item(16, 52).Value = 2
$ws.Cells.Item(17, 52).Value = 2
$ws.Cells.Item(23, 52).Value = 7
$ws.Cells.Item(25, 52).Value = 30
$ws.Cells.Item(26, 52).Value = 3
$ws.Cells.Item(28, 52).Value = 336
$ws.Cells.Item(29, 52).Value = 0
$ws.Cells.Item(31, 52).Value = 2
$ws.Cells.Item(35, 52).Value = 16
$ws.Cells.Item(36, 51).Value = 1
$ws.Cells.Item(36, 52).Value = 4
$ws.Cells.Item(38, 52).Value = 69
$ws.Cells.Item(41, 52).Value = 9
$ws.Cells.Item(42, 52).Value = 43
$ws.Cells.Item(44, 52).Value = 2
$ws.Cells.Item(45, 52).Value = 344
$ws.Cells.Item(46, 52).Value = 84
$ws.Cells.Item(47, 52).Value = 109
$ws.Cells.Item(48, 52).Value = 7
$ws.Cells.Item(49, 52).Value = 58
$ws.Cells.Item(50, 52).Value = 3
$ws.Cells.Item(51, 52).Value = 0
$ws.Cells.Item(54, 52).Value = 14
$ws.Cells.Item(55, 52).Value = 0
$ws.Cells.Item(56, 52).Value = 0
$ws.Cells.Item(58, 52).Value = 22
$ws.Cells.Item(59, 52).Value = 2

# --- Row 37: revised weekly counts for weeks 37-44 (columns AN-AU), plus
#     the new week-49 value.
$ws.Cells.Item(37, 40).Value = 5
$ws.Cells.Item(37, 41).Value = 10
$ws.Cells.Item(37, 43).Value = 11
$ws.Cells.Item(37, 44).Value = 4
$ws.Cells.Item(37, 45).Value = 10
$ws.Cells.Item(37, 46).Value = 10
$ws.Cells.Item(37, 47).Value = 11
$ws.Cells.Item(37, 52).Value = 4

# --- Row 43: revised weekly counts for weeks 15-31 (columns R-AH), plus
#     the new week-49 value.
$ws.Cells.Item(43, 18).Value = 61
$ws.Cells.Item(43, 19).Value = 35
$ws.Cells.Item(43, 20).Value = 45
$ws.Cells.Item(43, 21).Value = 40
$ws.Cells.Item(43, 22).Value = 41
$ws.Cells.Item(43, 23).Value = 48
$ws.Cells.Item(43, 24).Value = 70
$ws.Cells.Item(43, 25).Value = 98
$ws.Cells.Item(43, 26).Value = 59
$ws.Cells.Item(43, 27).Value = 80
$ws.Cells.Item(43, 28).Value = 71
$ws.Cells.Item(43, 30).Value = 39
$ws.Cells.Item(43, 31).Value = 56
$ws.Cells.Item(43, 32).Value = 51
$ws.Cells.Item(43, 33).Value = 52
$ws.Cells.Item(43, 34).Value = 57
$ws.Cells.Item(43, 52).Value = 49

# --- Row 57: revised counts for weeks 1 and 4 (columns D, G, H), plus the
#     new week-49 value.
$ws.Cells.Item(57, 4).Value = 7
$ws.Cells.Item(57, 7).Value = 8
$ws.Cells.Item(57, 8).Value = 6
$ws.Cells.Item(57, 52).Value = 8
